# Restore C10 on the "Rules" sheet to 1 (was 18)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
